# Workbook was uploaded with an updated "as of" date in the confidential
# disclaimer footnote plus refreshed Weight / Percent Change figures for
# the first few sector rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; temporarily unprotect so the cells can be
# written, then restore protection afterwards.
$ws.Unprotect()

# --- Update the "as of" date in the confidential disclaimer text (A9) ---
$cellA9 = $ws.Range("A9")
$currentText = $cellA9.Value2()
$updatedText = $currentText -replace "2021-05-17", "2021-05-18"
$cellA9.Value = $updatedText

# Re-entering a multi-line value can trigger an automatic row-height
# adjustment; AutoFit the row back so height metadata stays the way the
# workbook originally had it (auto/default row height, no explicit override).
$ws.Rows.Item(9).AutoFit()

# --- Update Weight (D) / Percent Change (E) figures for rows 2-6 ---
$ws.Range("D2").Value = 0.2586215167488806
$ws.Range("E2").Value = -0.01105096977898057

$ws.Range("D3").Value = 0.2554261369643501
$ws.Range("E3").Value = -0.013965744400527

$ws.Range("D4").Value = 0.2441173399592604
$ws.Range("E4").Value = -0.01455799252945122

$ws.Range("D5").Value = 0.2418350063275089
$ws.Range("E5").Value = -0.01019874476987448

$ws.Range("E6").Value = -0.0124455066252358

# Restore sheet protection (same options the workbook originally shipped
# with: objects/scenarios locked, but column/row formatting allowed).
$ws.Protect($null, $true, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $true, $true)
